# chore: update Sheets via scheduled runner
# Refreshes market/profit figures (currentAveragePrice*, LevePrice*, LeveProfit*)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4190.2
$ws.Range("J19").Value = 350
$ws.Range("L19").Value = 350
$ws.Range("N19").Value = -700
$ws.Range("H32").Value = 519.8
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 519.8
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 519.8
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -1171.8
$ws.Range("H93").Value = 34433.5
$ws.Range("J93").Value = 34433.5
$ws.Range("L93").Value = 34433.5
$ws.Range("N93").Value = -39425.5
$ws.Range("H129").Value = 823.5663
$ws.Range("I129").Value = 655.4
$ws.Range("J129").Value = 834.3461
$ws.Range("K129").Value = 1966.2
$ws.Range("L129").Value = 2503.0383
$ws.Range("M129").Value = 3033.8
$ws.Range("N129").Value = -12503.0383
$ws.Range("H138").Value = 2246.2542
$ws.Range("I138").Value = 1161.5172
$ws.Range("J138").Value = 3294.8333
$ws.Range("K138").Value = 3484.5516
$ws.Range("L138").Value = 9884.499899999999
$ws.Range("M138").Value = 1655.4484
$ws.Range("N138").Value = -20164.4999
$ws.Range("H141").Value = 2631.5
$ws.Range("I141").Value = 2698.889
$ws.Range("J141").Value = 2429.3333
$ws.Range("K141").Value = 8096.667
$ws.Range("L141").Value = 7287.999899999999
$ws.Range("M141").Value = -2916.667
$ws.Range("N141").Value = -17647.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29816.684
$ws.Range("I32").Value = 32088.854
$ws.Range("J32").Value = 10503.25
$ws.Range("K32").Value = 32088.854
$ws.Range("L32").Value = 10503.25
$ws.Range("M32").Value = -31801.854
$ws.Range("N32").Value = -11077.25
$ws.Range("H61").Value = 2828.6538
$ws.Range("J61").Value = 5999.5713
$ws.Range("L61").Value = 5999.5713
$ws.Range("N61").Value = -6423.5713
$ws.Range("H132").Value = 23381.959
$ws.Range("I132").Value = 1881.6923
$ws.Range("J132").Value = 48791.363
$ws.Range("K132").Value = 5645.0769
$ws.Range("L132").Value = 146374.089
$ws.Range("M132").Value = -3115.0769
$ws.Range("N132").Value = -151434.089
$ws.Range("H136").Value = 2828.6538
$ws.Range("J136").Value = 5999.5713
$ws.Range("L136").Value = 17998.7139
$ws.Range("N136").Value = -23098.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1749.9286
$ws.Range("I86").Value = 1584.1578
$ws.Range("J86").Value = 2099.889
$ws.Range("K86").Value = 1584.1578
$ws.Range("L86").Value = 2099.889
$ws.Range("M86").Value = -461.1578
$ws.Range("N86").Value = -4345.889
$ws.Range("H89").Value = 1749.9286
$ws.Range("I89").Value = 1584.1578
$ws.Range("J89").Value = 2099.889
$ws.Range("K89").Value = 7920.789
$ws.Range("L89").Value = 10499.445
$ws.Range("M89").Value = -2304.789
$ws.Range("N89").Value = -21731.445
$ws.Range("H134").Value = 62767
$ws.Range("I134").Value = 75900.86
$ws.Range("J134").Value = 1475.6666
$ws.Range("K134").Value = 227702.58
$ws.Range("L134").Value = 4426.9998
$ws.Range("M134").Value = -225167.58
$ws.Range("N134").Value = -9496.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 944.5172
$ws.Range("I16").Value = 714.55
$ws.Range("K16").Value = 714.55
$ws.Range("M16").Value = -427.55
$ws.Range("H31").Value = 12499.323
$ws.Range("I31").Value = 22542.715
$ws.Range("J31").Value = 5468.95
$ws.Range("K31").Value = 22542.715
$ws.Range("L31").Value = 5468.95
$ws.Range("M31").Value = -22247.715
$ws.Range("N31").Value = -6058.95
$ws.Range("H34").Value = 12499.323
$ws.Range("I34").Value = 22542.715
$ws.Range("J34").Value = 5468.95
$ws.Range("K34").Value = 22542.715
$ws.Range("L34").Value = 5468.95
$ws.Range("M34").Value = -22340.715
$ws.Range("N34").Value = -5872.95
$ws.Range("H58").Value = 21218.92
$ws.Range("I58").Value = 1301.5333
$ws.Range("J58").Value = 51095
$ws.Range("K58").Value = 1301.5333
$ws.Range("L58").Value = 51095
$ws.Range("M58").Value = -1098.5333
$ws.Range("N58").Value = -51501
$ws.Range("H113").Value = 944.5172
$ws.Range("I113").Value = 714.55
$ws.Range("K113").Value = 714.55
$ws.Range("M113").Value = 1455.45
$ws.Range("H136").Value = 21218.92
$ws.Range("I136").Value = 1301.5333
$ws.Range("J136").Value = 51095
$ws.Range("K136").Value = 3904.5999
$ws.Range("L136").Value = 153285
$ws.Range("M136").Value = -1354.5999
$ws.Range("N136").Value = -158385

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 3522.5
$ws.Range("J54").Value = 3522.5
$ws.Range("L54").Value = 10567.5
$ws.Range("N54").Value = -11685.5
$ws.Range("H122").Value = 774.7143
$ws.Range("I122").Value = 524.8
$ws.Range("J122").Value = 1399.5
$ws.Range("K122").Value = 4723.2
$ws.Range("L122").Value = 12595.5
$ws.Range("M122").Value = -2273.2
$ws.Range("N122").Value = -17495.5
$ws.Range("H131").Value = 764.62
$ws.Range("J131").Value = 764.62
$ws.Range("L131").Value = 2293.86
$ws.Range("N131").Value = -12373.86

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3717.8518
$ws.Range("I80").Value = 3118.5454
$ws.Range("J80").Value = 4129.875
$ws.Range("K80").Value = 3118.5454
$ws.Range("L80").Value = 4129.875
$ws.Range("M80").Value = -2120.5454
$ws.Range("N80").Value = -6125.875
$ws.Range("H83").Value = 3717.8518
$ws.Range("I83").Value = 3118.5454
$ws.Range("J83").Value = 4129.875
$ws.Range("K83").Value = 15592.727
$ws.Range("L83").Value = 20649.375
$ws.Range("M83").Value = -10600.727
$ws.Range("N83").Value = -30633.375
$ws.Range("H113").Value = 2227.1155
$ws.Range("I113").Value = 1662.6111
$ws.Range("K113").Value = 1662.6111
$ws.Range("M113").Value = 507.3888999999999
$ws.Range("I126").Value = 3221.875
$ws.Range("J126").Value = 4600.8237
$ws.Range("K126").Value = 9665.625
$ws.Range("L126").Value = 13802.4711
$ws.Range("M126").Value = -7195.625
$ws.Range("N126").Value = -18742.4711
$ws.Range("H136").Value = 24581.5
$ws.Range("J136").Value = 24581.5
$ws.Range("L136").Value = 73744.5
$ws.Range("N136").Value = -78844.5
$ws.Range("H141").Value = 41333
$ws.Range("J141").Value = 41333
$ws.Range("L141").Value = 41333
$ws.Range("N141").Value = -51693

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 5012500
$ws.Range("I38").Value = 20000000
$ws.Range("J38").Value = 16666.666
$ws.Range("K38").Value = 20000000
$ws.Range("L38").Value = 16666.666
$ws.Range("M38").Value = -19999590
$ws.Range("N38").Value = -17486.666
$ws.Range("H132").Value = 806724.9
$ws.Range("I132").Value = 1098079.5
$ws.Range("J132").Value = 5499.5
$ws.Range("K132").Value = 3294238.5
$ws.Range("L132").Value = 16498.5
$ws.Range("M132").Value = -3291708.5
$ws.Range("N132").Value = -21558.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1004.8571
$ws.Range("J81").Value = 461.14285
$ws.Range("L81").Value = 922.2857
$ws.Range("N81").Value = -3044.2857
$ws.Range("H84").Value = 1004.8571
$ws.Range("J84").Value = 461.14285
$ws.Range("L84").Value = 4611.4285
$ws.Range("N84").Value = -15219.4285
$ws.Range("H107").Value = 775.16
$ws.Range("I107").Value = 314
$ws.Range("J107").Value = 1274.75
$ws.Range("K107").Value = 942
$ws.Range("L107").Value = 3824.25
$ws.Range("M107").Value = 978
$ws.Range("N107").Value = -7664.25
$ws.Range("H122").Value = 1570.96
$ws.Range("I122").Value = 1391.0526
$ws.Range("J122").Value = 2140.6667
$ws.Range("K122").Value = 4173.1578
$ws.Range("L122").Value = 6422.000100000001
$ws.Range("M122").Value = -1723.1578
$ws.Range("N122").Value = -11322.0001
